$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 -- append two spaces to the existing sentence,
# then add three colored runs that together read:
# "(This is a change - Version for branch alternate)"
$p1 = $d.Paragraphs(1).Range
$end = $p1.End - 1
$r = $d.Range($end, $end)
$r.InsertAfter("  ")

$t1 = "(This is a change " + [char]0x2013 + " Ve"
$t2 = "rsion for branch alternate"
$t3 = ")"

$pos = $end + 2
$r = $d.Range($pos, $pos)
$r.InsertAfter($t1)
$r1 = $d.Range($pos, $pos + $t1.Length)
$r1.Font.Color = 192

$pos = $pos + $t1.Length
$r = $d.Range($pos, $pos)
$r.InsertAfter($t2)
$r2 = $d.Range($pos, $pos + $t2.Length)
$r2.Font.Color = 192

$pos = $pos + $t2.Length
$r = $d.Range($pos, $pos)
$r.InsertAfter($t3)
$r3 = $d.Range($pos, $pos + $t3.Length)
$r3.Font.Color = 192

# --- Edit 2: insert a new, empty, shaded paragraph right after the
# "It will be treated as a binary file by Git." paragraph.
$p2 = $d.Paragraphs(2).Range
$p2.InsertParagraphAfter()
$newPara = $d.Paragraphs(3)
$newPara.Range.Font.Name = "Calibri"
$newPara.Range.Font.Bold = $true
$newPara.Range.Font.Color = 2236704
$newPara.Range.Font.NameFarEast = "Times New Roman"
$newPara.Range.Font.NameBi = "Calibri"
$newPara.Range.Font.BoldBi = $true
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457
